# The deck ships two themes:
#   theme1.xml -> "Office Theme" (Office color scheme) - used only by the Notes Master
#   theme2.xml -> "Integral"     (Red Violet color scheme) - used by the Slide Master
#                                  (and therefore by every slide in the deck)
#
# The authored change swaps the content of theme1.xml and theme2.xml, which in
# effect re-colors the actual slide deck (theme2.xml, reachable via
# SlideMaster.ColorScheme) from the "Red Violet" palette to the "Office" palette.
#
# PowerPoint's ColorScheme.Colors(n).RGB uses the standard OLE RGB encoding
# (R + G*256 + B*65536), and the 12 slots map 1:1 onto the theme's <a:clrScheme>
# children in document order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
#
# Target values below are the "Office" theme colors (000000, FFFFFF, 44546A,
# E7E6E6, 5B9BD5, ED7D31, A5A5A5, FFC000, 4472C4, 70AD47, 0563C1, 954F72).

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72

# Best-effort: restore the scheme/theme display names to match the "Office"
# theme (harmless if the host treats these as read-only for serialization).
$cs.Name = "Office"
